$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.787.23"
$ws.Range("E2").Value = "  +0.93%  "

$ws.Range("D3").Value = "2.505.19"
$ws.Range("E3").Value = "  +0.04%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.90"
$ws.Range("E5").Value = "  -0.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.79"
$ws.Range("E6").Value = "  -1.50%  "

$ws.Range("E7").Value = "  -0.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.558"
$ws.Range("E9").Value = "  +2.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.21"
$ws.Range("E10").Value = "  +3.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.64"
$ws.Range("E11").Value = "  +5.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0813"
$ws.Range("E12").Value = "  -0.44%  "

$ws.Range("E13").Value = "  +0.66%  "

$ws.Range("E14").Value = "  -0.81%  "

$ws.Range("D15").Value = "2.899.88"
$ws.Range("E15").Value = "  +0.16%  "

$ws.Range("D16").Value = "2.502.63"
$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.849"
$ws.Range("E17").Value = "  -1.17%  "

$ws.Range("D18").Value = "47.690.31"
$ws.Range("E18").Value = "  +0.82%  "

$ws.Range("E19").Value = "  +2.60%  "

$ws.Range("E20").Value = "  -1.77%  "

$ws.Range("D21").Value = "0.0₃0940"
$ws.Range("E21").Value = "  -0.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.76"
$ws.Range("E22").Value = "  +7.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.90"
$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "247.11"
$ws.Range("E24").Value = "  -1.08%  "

$ws.Range("E25").Value = "  -0.65%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.75"
$ws.Range("E27").Value = "  -1.53%  "

$ws.Range("E30").Value = "  +4.59%  "

$ws.Range("E31").Value = "  -2.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.81"
$ws.Range("E32").Value = "  -0.59%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.00"
$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("E36").Value = "  -1.28%  "

$ws.Range("E37").Value = "  -1.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.67"
$ws.Range("E38").Value = "  -1.27%  "

$ws.Range("E39").Value = "  -1.13%  "

$ws.Range("E40").Value = "  -0.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.28"
$ws.Range("E41").Value = "  +4.40%  "

$ws.Range("E42").Value = "  -2.40%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "118.59"
$ws.Range("E43").Value = "  -2.79%  "

$ws.Range("E44").Value = "  -0.68%  "

$ws.Range("D45").Value = "1.999.38"
$ws.Range("E45").Value = "  +0.42%  "

$ws.Range("E46").Value = "  +1.30%  "

$ws.Range("E47").Value = "  -2.98%  "

$ws.Range("E48").Value = "  +0.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.09"
$ws.Range("E49").Value = "  +0.29%  "

$ws.Range("E50").Value = "  -2.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.48"
$ws.Range("E51").Value = "  +1.81%  "

# Row 28/29: Cosmos/Toncoin swap
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.31"
$ws.Range("E28").Value = "  +0.67%  "

$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.24"
$ws.Range("E29").Value = "  +2.02%  "
